$d = $word.ActiveDocument

function Scope-Range([string]$anchorText) {
    # Find a unique anchor string and return a Range bounded exactly to it
    # (Start..End), positioned so a subsequent Find on it only matches
    # within those bounds.
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText)
    if (-not $found) {
        throw "anchor not found: $anchorText"
    }
    $s = $rng.Start
    $e = $rng.End
    $rng.Collapse(1)
    $rng.MoveEnd(1, ($e - $s))
    return $rng
}

function Split-Off([string]$anchorText, [string]$secondPartText) {
    # Splits the run(s) covering $anchorText so that $secondPartText
    # (a trailing substring of $anchorText) becomes its own run, by
    # toggling a direct-character-formatting property on/off (which
    # forces the run boundary without altering the visible formatting).
    $scoped = Scope-Range $anchorText
    $found = $scoped.Find.Execute($secondPartText)
    if (-not $found) {
        throw "second part not found within anchor: $secondPartText"
    }
    $scoped.Bold = $true
    $scoped.Bold = $false
}

# --- word/document.xml text edits -----------------------------------

# Paragraph: "...output/cohort-stats/{species}-{statistic}-{timestep}.gis"
# Split "/cohort-stats/{species}-{statistic}-{" into
#   "/cohort-stats/{spe" + "cies}-{statistic}-{"
Split-Off "/cohort-stats/{species}-{statistic}-{timestep}.gis" "cies}-{statistic}-{"

# Paragraph: "   output/cohort-stats/AGE-{statistic}-{timestep}.gis"
# Split "   output/cohort-stats/AGE-{statistic}-{" into
#   "   output/cohort-stat" + "s/AGE-{statistic}-{"
Split-Off "   output/cohort-stats/AGE-{statistic}-{timestep}.gis" "s/AGE-{statistic}-{"

# Paragraph: "   output/cohort-stats/SPP-{statistic}-{timestep}.gis"
# Split "timestep" into "time" + "step"
Split-Off "   output/cohort-stats/SPP-{statistic}-{timestep}.gis" "step"

# All three "gis" (file-extension) runs become "img".
$d.Content.Find.Execute("gis", $true, $false, $false, $false, $false, $true, 1, $false, "img", 2) | Out-Null

# --- word/styles.xml: add <w:rsid w:val="00297BD4"/> to (almost) every
#     style, and tweak DefaultParagraphFont's uiPriority/unhideWhenUsed --

$rsidVal = "00297BD4"
$styles = $d.Styles
for ($i = 1; $i -le $styles.Count; $i++) {
    $st = $styles.Item($i)
    if ($st.BuiltIn -and $st.NameLocal -eq "Default Paragraph Font") {
        continue
    }
    try {
        $st.AddRsid($rsidVal)
    } catch {
    }
}

Write-Output "done"
